# Auto-generated: apply scheduled-runner market-data updates to Sheets/Seraph_Profits.xlsx
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds Leve crafting-profit data;
# columns H-N are live market prices / profit calcs refreshed by the runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 42449
$ws.Range("J51").Value = 42449
$ws.Range("L51").Value = 42449
$ws.Range("N51").Value = -43417
$ws.Range("H54").Value = 3000
$ws.Range("I54").Value = 3000
$ws.Range("K54").Value = 3000
$ws.Range("M54").Value = -2514
$ws.Range("H70").Value = 5122.0884
$ws.Range("I70").Value = 6070.9473
$ws.Range("J70").Value = 3920.2
$ws.Range("K70").Value = 18212.8419
$ws.Range("L70").Value = 11760.6
$ws.Range("M70").Value = -17942.8419
$ws.Range("N70").Value = -12300.6
$ws.Range("H73").Value = 5122.0884
$ws.Range("I73").Value = 6070.9473
$ws.Range("J73").Value = 3920.2
$ws.Range("K73").Value = 18212.8419
$ws.Range("L73").Value = 11760.6
$ws.Range("M73").Value = -17276.8419
$ws.Range("N73").Value = -13632.6
$ws.Range("H98").Value = 715.06665
$ws.Range("I98").Value = 751.8570999999999
$ws.Range("J98").Value = 200
$ws.Range("K98").Value = 751.8570999999999
$ws.Range("L98").Value = 200
$ws.Range("M98").Value = 746.1429000000001
$ws.Range("N98").Value = -3196
$ws.Range("H103").Value = 1636.875
$ws.Range("J103").Value = 1478.5
$ws.Range("L103").Value = 4435.5
$ws.Range("N103").Value = -5607.5
$ws.Range("H122").Value = 715.06665
$ws.Range("I122").Value = 751.8570999999999
$ws.Range("J122").Value = 200
$ws.Range("K122").Value = 2255.5713
$ws.Range("L122").Value = 600
$ws.Range("M122").Value = 194.4287000000004
$ws.Range("N122").Value = -5500
$ws.Range("H125").Value = 7110
$ws.Range("J125").Value = 7181
$ws.Range("L125").Value = 64629
$ws.Range("N125").Value = -69549
$ws.Range("H131").Value = 3650.2632
$ws.Range("I131").Value = 422.36365
$ws.Range("J131").Value = 8088.625
$ws.Range("K131").Value = 1267.09095
$ws.Range("L131").Value = 24265.875
$ws.Range("M131").Value = 3772.90905
$ws.Range("N131").Value = -34345.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H97").Value = 469.69232
$ws.Range("I97").Value = 469.69232
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 469.69232
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 26.30768
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 5999
$ws.Range("I102").Value = 5999
$ws.Range("K102").Value = 5999
$ws.Range("M102").Value = -4377
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 3759.35
$ws.Range("I122").Value = 3129.8462
$ws.Range("K122").Value = 9389.5386
$ws.Range("M122").Value = -6939.5386
$ws.Range("H132").Value = 1920.1111
$ws.Range("I132").Value = 1396.5333
$ws.Range("K132").Value = 4189.5999
$ws.Range("M132").Value = -1659.5999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1694.6316
$ws.Range("I99").Value = 1386.5333
$ws.Range("K99").Value = 1386.5333
$ws.Range("M99").Value = 111.4666999999999
$ws.Range("H105").Value = 3809.4814
$ws.Range("I105").Value = 2921.9524
$ws.Range("J105").Value = 6915.8335
$ws.Range("K105").Value = 2921.9524
$ws.Range("L105").Value = 6915.8335
$ws.Range("M105").Value = -1174.9524
$ws.Range("N105").Value = -10409.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2411
$ws.Range("I16").Value = 2610.5
$ws.Range("J16").Value = 2012
$ws.Range("K16").Value = 2610.5
$ws.Range("L16").Value = 2012
$ws.Range("M16").Value = -2323.5
$ws.Range("N16").Value = -2586
$ws.Range("H58").Value = 2485
$ws.Range("I58").Value = 1919.4736
$ws.Range("K58").Value = 1919.4736
$ws.Range("M58").Value = -1716.4736
$ws.Range("H95").Value = 16151.375
$ws.Range("J95").Value = 16151.375
$ws.Range("L95").Value = 16151.375
$ws.Range("N95").Value = -21643.375
$ws.Range("H107").Value = 990.119
$ws.Range("I107").Value = 729.6667
$ws.Range("K107").Value = 729.6667
$ws.Range("M107").Value = 1190.3333
$ws.Range("H113").Value = 2411
$ws.Range("I113").Value = 2610.5
$ws.Range("J113").Value = 2012
$ws.Range("K113").Value = 2610.5
$ws.Range("L113").Value = 2012
$ws.Range("M113").Value = -440.5
$ws.Range("N113").Value = -6352
$ws.Range("H122").Value = 1685.4546
$ws.Range("I122").Value = 1521.5714
$ws.Range("J122").Value = 1972.25
$ws.Range("K122").Value = 4564.7142
$ws.Range("L122").Value = 5916.75
$ws.Range("M122").Value = -2114.7142
$ws.Range("N122").Value = -10816.75
$ws.Range("H134").Value = 1494.8889
$ws.Range("I134").Value = 1079.6154
$ws.Range("K134").Value = 3238.8462
$ws.Range("M134").Value = -703.8462
$ws.Range("H136").Value = 2485
$ws.Range("I136").Value = 1919.4736
$ws.Range("K136").Value = 5758.4208
$ws.Range("M136").Value = -3208.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 146
$ws.Range("I26").Value = 26.666666
$ws.Range("J26").Value = 217.6
$ws.Range("K26").Value = 79.99999800000001
$ws.Range("L26").Value = 652.8
$ws.Range("M26").Value = 208.000002
$ws.Range("N26").Value = -1228.8
$ws.Range("H107").Value = 111391.22
$ws.Range("I107").Value = 310.33334
$ws.Range("J107").Value = 166931.67
$ws.Range("K107").Value = 931.0000200000001
$ws.Range("L107").Value = 500795.01
$ws.Range("M107").Value = 988.9999799999999
$ws.Range("N107").Value = -504635.01
$ws.Range("H132").Value = 1974.5
$ws.Range("I132").Value = 1974.5
$ws.Range("K132").Value = 17770.5
$ws.Range("M132").Value = -15240.5
$ws.Range("H133").Value = 4312.25
$ws.Range("I133").Value = 749.6667
$ws.Range("J133").Value = 15000
$ws.Range("K133").Value = 2249.0001
$ws.Range("L133").Value = 45000
$ws.Range("M133").Value = 2810.9999
$ws.Range("N133").Value = -55120
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2004.4783
$ws.Range("I102").Value = 546.86664
$ws.Range("J102").Value = 4737.5
$ws.Range("K102").Value = 546.86664
$ws.Range("L102").Value = 4737.5
$ws.Range("M102").Value = 1075.13336
$ws.Range("N102").Value = -7981.5
$ws.Range("H107").Value = 1649.95
$ws.Range("I107").Value = 1405.8235
$ws.Range("K107").Value = 1405.8235
$ws.Range("M107").Value = 514.1765
$ws.Range("H126").Value = 4444.727
$ws.Range("I126").Value = 3799.6
$ws.Range("J126").Value = 4982.3335
$ws.Range("K126").Value = 11398.8
$ws.Range("L126").Value = 14947.0005
$ws.Range("M126").Value = -8928.799999999999
$ws.Range("N126").Value = -19887.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2187.5
$ws.Range("I46").Value = 1250
$ws.Range("K46").Value = 1250
$ws.Range("M46").Value = -1062
$ws.Range("H93").Value = 4110.5557
$ws.Range("I93").Value = 4249.375
$ws.Range("K93").Value = 4249.375
$ws.Range("M93").Value = -3001.375
$ws.Range("H100").Value = 8571.286
$ws.Range("I100").Value = 5833.5
$ws.Range("K100").Value = 5833.5
$ws.Range("M100").Value = -5292.5
$ws.Range("H136").Value = 6906.3076
$ws.Range("I136").Value = 5311.75
$ws.Range("K136").Value = 15935.25
$ws.Range("M136").Value = -13385.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3743.2727
$ws.Range("I132").Value = 3691.889
$ws.Range("J132").Value = 3974.5
$ws.Range("K132").Value = 11075.667
$ws.Range("L132").Value = 11923.5
$ws.Range("M132").Value = -8545.667000000001
$ws.Range("N132").Value = -16983.5
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 8106.875
$ws.Range("I136").Value = 8978.571
$ws.Range("K136").Value = 26935.713
$ws.Range("M136").Value = -24385.713
